$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format D and E columns (and B/C for the swapped rows) as Text so that
# numeric-looking strings (e.g. "0.9997", "0.000008941") are preserved exactly
# as text instead of being auto-converted to numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('B48').NumberFormat = '@'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('B50').NumberFormat = '@'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '29.182.31'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.827.29'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '241.86'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').Value = '0.6205'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.07332'
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('D9').Value = '0.2902'
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('D10').Value = '23.04'
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('D11').Value = '0.07691'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '1.825.99'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = '0.6650'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = '82.36'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').Value = '0.000008941'
$ws.Range('E16').Value = '  -4.63%  '
$ws.Range('D17').Value = '5.861'
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('D18').Value = '29.162.65'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').Value = '2.073.17'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').Value = '238.26'
$ws.Range('E20').Value = '  +6.87%  '
$ws.Range('D21').Value = '12.44'
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '7.310'
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '158.38'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').Value = '0.1424'
$ws.Range('E26').Value = '  +2.32%  '
$ws.Range('D27').Value = '8.480'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = '17.68'
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').Value = '1.485'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').Value = '0.05586'
$ws.Range('E30').Value = '  -3.13%  '
$ws.Range('D31').Value = '4.088'
$ws.Range('E31').Value = '  -0.83%  '
$ws.Range('D32').Value = '4.097'
$ws.Range('E32').Value = '  -1.32%  '
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('D34').Value = '1.844'
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').Value = '0.7351'
$ws.Range('E35').Value = '  -0.33%  '
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('D37').Value = '2.628'
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('D38').Value = '2.842'
$ws.Range('E38').Value = '  +2.91%  '
$ws.Range('D39').Value = '1.211.64'
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('D40').Value = '0.01768'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').Value = '6.307'
$ws.Range('E41').Value = '  -2.66%  '
$ws.Range('D42').Value = '0.9208'
$ws.Range('E42').Value = '  +3.40%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').Value = '101.68'
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('D45').Value = '1.978.14'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').Value = '64.86'
$ws.Range('E46').Value = '  -1.57%  '
$ws.Range('D47').Value = '0.5091'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.00000000118'
$ws.Range('E48').Value = '  -5.16%  '
$ws.Range('D49').Value = '0.4027'
$ws.Range('E49').Value = '  -0.60%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '9.115'
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('D51').Value = '0.05761'
$ws.Range('E51').Value = '  -1.13%  '
